# The commit adds one new weekly price observation for "Jengibre" at the
# top of the historical (date-descending) data block, pushing every
# existing record from row 57 downward by one row (old row 57 -> new row
# 58, ..., old row 136 -> new row 137).
#
# Reproduce that with a real row insert (so formatting/styles on the
# shifted rows carry down naturally), then populate the freshly inserted
# row 57 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 57; rows 57:136 shift to 58:137.
$ws.Rows.Item(57).Insert()

# Fill in the new record in row 57.
$ws.Range("A57").Value = 6
$ws.Range("B57").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C57").Value = "Metropolitana"
$ws.Range("D57").Value = 45117
$ws.Range("E57").Value = 13
$ws.Range("F57").Value = 100114007
$ws.Range("G57").Value = "Jengibre"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 300
$ws.Range("K57").Value = 13000
$ws.Range("L57").Value = 14000
$ws.Range("M57").Value = 13567
$ws.Range("N57").Value = "`$/caja 13 kilos"
$ws.Range("O57").Value = "Per" + [char]0x00FA
$ws.Range("P57").Value = 1044
$ws.Range("Q57").Value = 13
$ws.Range("R57").Value = "Hortaliza"
